$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.094.26"
$ws.Range("E2").Value = "  +0.48%  "

$ws.Range("D3").Value = "2.982.89"
$ws.Range("E3").Value = "  +1.71%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "353.80"
$ws.Range("E5").Value = "  +0.34%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.18"
$ws.Range("E6").Value = "  -3.55%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.559"
$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.616"
$ws.Range("E9").Value = "  -0.74%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.37"
$ws.Range("E10").Value = "  -2.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.139"
$ws.Range("E11").Value = "  +1.40%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0854"
$ws.Range("E12").Value = "  -4.04%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.19"
$ws.Range("E13").Value = "  -4.27%  "

$ws.Range("D14").Value = "3.439.16"
$ws.Range("E14").Value = "  +1.32%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.62"
$ws.Range("E15").Value = "  -1.74%  "

$ws.Range("D16").Value = "2.977.09"
$ws.Range("E16").Value = "  +1.45%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.00"
$ws.Range("E17").Value = "  +1.52%  "

$ws.Range("D18").Value = "52.071.44"
$ws.Range("E18").Value = "  +0.30%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.46"
$ws.Range("E19").Value = "  +4.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.47"
$ws.Range("E20").Value = "  -1.85%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.57"
$ws.Range("E21").Value = "  -4.55%  "

$ws.Range("D22").Value = "0.0₃0973"
$ws.Range("E22").Value = "  -1.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.53"
$ws.Range("E23").Value = "  -2.42%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "264.32"
$ws.Range("E24").Value = "  -1.58%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.73"
$ws.Range("E25").Value = "  -1.91%  "

$ws.Range("E26").Value = "  -2.51%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.84"
$ws.Range("E27").Value = "  -0.30%  "

$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.48"
$ws.Range("E28").Value = "  +2.87%  "

$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.107"
$ws.Range("E30").Value = "  +2.55%  "

$ws.Range("E31").Value = "  -2.55%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.08"
$ws.Range("E32").Value = "  -3.46%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "36.22"
$ws.Range("E33").Value = "  -2.52%  "

$ws.Range("E34").Value = "  -4.76%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.65"
$ws.Range("E35").Value = "  -4.40%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0437"
$ws.Range("E36").Value = "  -3.92%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.997"
$ws.Range("E37").Value = "  -0.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.21"
$ws.Range("E38").Value = "  -3.92%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.87"
$ws.Range("E39").Value = "  -4.78%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.97"
$ws.Range("E40").Value = "  -4.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.72"
$ws.Range("E41").Value = "  +1.19%  "

$ws.Range("E42").Value = "  -0.97%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "123.19"
$ws.Range("E43").Value = "  +10.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.51"
$ws.Range("E44").Value = "  -2.85%  "

$ws.Range("E45").Value = "  -3.74%  "

$ws.Range("D46").Value = "2.114.88"
$ws.Range("E46").Value = "  -2.77%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.35"
$ws.Range("E47").Value = "  -4.46%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.31"
$ws.Range("E48").Value = "  -8.83%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.239"
$ws.Range("E49").Value = "  -4.19%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0342"
$ws.Range("E50").Value = "  -2.05%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.928"
$ws.Range("E51").Value = "  -1.66%  "
